$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new bullet ("Subsidy causes change in survival rates")
#    right before the "Ceteris Paribus (Business conditions)" bullet,
#    at the same outline level (ListParagraph / ilvl 1 / numId 2).
# ------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Ceteris Paribus") | Out-Null
$ceterisPara = $findRng.Paragraphs(1)
[void]$ceterisPara.Range.InsertParagraphBefore()

# Re-locate the "Ceteris Paribus" paragraph (indices shifted after insert)
# so we can grab the new, still-empty paragraph right before it.
$findRng2 = $d.Content
$findRng2.Find.Execute("Ceteris Paribus") | Out-Null
$ceterisPara2 = $findRng2.Paragraphs(1)
$subsidyPara = $ceterisPara2.Previous()

$subsidyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">ubsidy </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>causes change in survival rates</w:t></w:r></w:p>
'@
$subsidyPara.Range.InsertXML($subsidyXml) | Out-Null

# ------------------------------------------------------------------
# 2) Rewrite the "Ceteris Paribus (Business conditions)" bullet's
#    text to "All businesses in equal situations" and add the
#    _GoBack bookmark at its end.
# ------------------------------------------------------------------
$findRng3 = $d.Content
$findRng3.Find.Execute("Ceteris Paribus") | Out-Null
$ceterisPara3 = $findRng3.Paragraphs(1)

$ceterisXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve">All businesses </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve">in </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t>equal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:bCs/></w:rPr><w:t xml:space="preserve"> situations</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$ceterisPara3.Range.InsertXML($ceterisXml) | Out-Null

# ------------------------------------------------------------------
# 3) The (previously empty) paragraph right after the "Ceteris
#    Paribus" bullet now holds a single bold space character.
# ------------------------------------------------------------------
$findRng4 = $d.Content
$findRng4.Find.Execute("All businesses in equal situations") | Out-Null
$ceterisPara4 = $findRng4.Paragraphs(1)
$afterPara = $ceterisPara4.Next()

$spaceXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$afterPara.Range.InsertXML($spaceXml) | Out-Null

# ------------------------------------------------------------------
# 4) Collapse the split runs in the two "Calculations for ..." lines
#    into single runs (text unchanged, but re-"typing" it over itself
#    via Find/Replace merges the runs into one).
# ------------------------------------------------------------------
$d.Content.Find.Execute("Calculations for subsidized clinics below. Calculation in attached excel file.", $true, $false, $false, $false, $false, $true, 1, $false, "Calculations for subsidized clinics below. Calculation in attached excel file.", 2) | Out-Null
$d.Content.Find.Execute("Calculations for ROI below. Calculation in attached excel file.", $true, $false, $false, $false, $false, $true, 1, $false, "Calculations for ROI below. Calculation in attached excel file.", 2) | Out-Null

$d.Save()
